$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PruebasUnit")
$ws.Range("A1").Value = "test"
